# LOQ4235.xlsx syllabus update
#
# Summary of the change (see commit diff):
#  - Row "Objetivos:" loses its long paragraph; B/C now hold the
#    professor's name/id instead.
#  - The row that used to hold only the professor's name/id is dropped
#    (its content is now reused one row up), so every label below it
#    shifts up by one row.
#  - "Programa resumido:" now pairs with "Semestral" instead of the old
#    free-text paragraph.
#  - The long "Programa:" paragraph is dropped too (content reused one
#    row up again - "01/01/2015" ends up as the Programa: value, matching
#    the source diff literally).
#  - A new "Critério:" row is introduced holding the
#    "O desenvolvimento..." paragraph (which used to sit on "Método:").
#  - "Método:" now pairs with the professor's name/id.
#  - The very last row ("Bibliografia:" / long bibliography paragraph)
#    is removed outright, which is what shrinks the sheet from 22 to 21
#    rows; "Bibliografia:" is reused as the new label for the
#    "Prova única..." text.
#
# Net effect on row count: 22 -> 21 rows (dimension A1:C22 -> A1:C21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the very last row outright - this is what actually shrinks the
# sheet. Everything from row 13 down shifts up by one.
$ws.Rows.Item(22).Delete()

# ---- Row 10 ("Objetivos:") ------------------------------------------------
# B/C: long paragraph -> professor name/id (cells already existed, so the
# existing style (s="2"/"s=3") is preserved automatically).
$ws.Cells.Item(10, 2).Value2 = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Cells.Item(10, 3).Value2 = "5840560 - Marco Antonio Carvalho Pereira"

# ---- Row 13 ("Programa resumido:") ----------------------------------------
# A13 is a brand-new cell (row 13 used to only have B/C). Copy the look of
# the other label cells in column A onto it before setting its value, so it
# lands on the same style record (s="1") instead of the engine's default.
$ws.Cells.Item(10, 1).Copy()
$ws.Cells.Item(13, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Cells.Item(13, 1).Value2 = "Programa resumido:"
$ws.Cells.Item(13, 2).Value2 = "Semestral"
$ws.Cells.Item(13, 3).Value2 = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# ---- Row 14 ("Short syllabus:") -------------------------------------------
$ws.Cells.Item(14, 1).Value2 = "Short syllabus:"
$ws.Cells.Item(14, 2).Clear()
$ws.Cells.Item(14, 3).Clear()

# ---- Row 15 ("Programa:") --------------------------------------------------
$ws.Cells.Item(15, 1).Value2 = "Programa:"
# B15/C15 are brand-new cells whose value ("01/01/2015") looks like a date,
# which Excel would otherwise auto-convert to a date serial number. Paste
# the VALUE from an existing text cell that already holds that exact string
# (row 8) so it is carried over as plain text/shared-string, then paste the
# FORMAT (only) from row 10 so the cell lands on the normal paragraph style
# instead of the engine's default-new-cell style.
$ws.Cells.Item(8, 2).Copy()
$ws.Cells.Item(15, 2).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(8, 3).Copy()
$ws.Cells.Item(15, 3).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(10, 2).Copy()
$ws.Cells.Item(15, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(10, 3).Copy()
$ws.Cells.Item(15, 3).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Rows.Item(15).RowHeight = 120

# ---- Row 16 ("Syllabus:") --------------------------------------------------
$ws.Cells.Item(16, 1).Value2 = "Syllabus:"
$ws.Cells.Item(16, 2).Clear()
$ws.Cells.Item(16, 3).Clear()

# ---- Row 17 ("Avaliação:") -------------------------------------------------
$ws.Cells.Item(17, 1).Value2 = "Avaliação:"
$ws.Rows.Item(17).AutoFit() | Out-Null

# ---- Row 18 ("Método:") -----------------------------------------------------
$ws.Cells.Item(18, 1).Value2 = "Método:"
# B18/C18 are brand-new cells - copy formats first, same trick as row 13/15.
$ws.Cells.Item(10, 2).Copy()
$ws.Cells.Item(18, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(10, 3).Copy()
$ws.Cells.Item(18, 3).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Cells.Item(18, 2).Value2 = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Cells.Item(18, 3).Value2 = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Rows.Item(18).RowHeight = 60

# ---- Row 19 ("Critério:") ---------------------------------------------------
$ws.Cells.Item(19, 1).Value2 = "Critério:"
$ws.Cells.Item(19, 2).Value2 = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."
$ws.Cells.Item(19, 3).Value2 = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."

# ---- Row 20 ("Norma de recuperação:") ---------------------------------------
$ws.Cells.Item(20, 1).Value2 = "Norma de recuperação:"
$ws.Cells.Item(20, 2).Value2 = "Provas e trabalhos."
$ws.Cells.Item(20, 3).Value2 = "Provas e trabalhos."

# ---- Row 21 ("Bibliografia:") -----------------------------------------------
$ws.Cells.Item(21, 1).Value2 = "Bibliografia:"
$ws.Cells.Item(21, 2).Value2 = "Prova única com nota maior ou igual a 5,0 (cinco)."
$ws.Cells.Item(21, 3).Value2 = "Prova única com nota maior ou igual a 5,0 (cinco)."
$ws.Rows.Item(21).RowHeight = 120
